$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first header from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Delete the old "Run 50" column (AZ). This shifts the old "Mean" column
# (BA) left into AZ, carrying its style/values with it.
$ws.Range("AZ:AZ").EntireColumn.Delete()

# Update the MaxFES (column A) values for rows 2-14
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Update the recomputed Mean column (now AZ, after the column delete) values
$ws.Range("AZ2").Value = 255.37499698
$ws.Range("AZ3").Value = 185.31951652
$ws.Range("AZ4").Value = 17.46175836
$ws.Range("AZ5").Value = 0.29901232
$ws.Range("AZ6").Value = 0.14622935
$ws.Range("AZ7").Value = 0.09632119
$ws.Range("AZ8").Value = 0.0720013
$ws.Range("AZ9").Value = 0.05755879
$ws.Range("AZ10").Value = 0.04624972
$ws.Range("AZ11").Value = 0.03830892
$ws.Range("AZ12").Value = 0.0343615
$ws.Range("AZ13").Value = 0.02985142
$ws.Range("AZ14").Value = 0.0260045
